$wb = $excel.ActiveWorkbook

# Helper: set a cell value by sheet name, row, column letter
function Set-Cell {
    param($ws, [string]$col, [int]$row, $value)
    $ws.Range("$col$row").Value = $value
}

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
Set-Cell $ws "H" 64 3916.6667
Set-Cell $ws "I" 64 3916.6667
Set-Cell $ws "J" 64 0
Set-Cell $ws "K" 64 3916.6667
Set-Cell $ws "L" 64 0
Set-Cell $ws "M" 64 -3668.6667

Set-Cell $ws "H" 67 3916.6667
Set-Cell $ws "I" 67 3916.6667
Set-Cell $ws "J" 67 0
Set-Cell $ws "K" 67 3916.6667
Set-Cell $ws "L" 67 0
Set-Cell $ws "M" 67 -3058.6667

Set-Cell $ws "H" 76 3246.1155
Set-Cell $ws "I" 76 3260.8262
Set-Cell $ws "J" 76 3133.3333
Set-Cell $ws "K" 76 3260.8262
Set-Cell $ws "L" 76 3133.3333
Set-Cell $ws "M" 76 -2945.8262
Set-Cell $ws "N" 76 -3763.3333

Set-Cell $ws "H" 79 3246.1155
Set-Cell $ws "I" 79 3260.8262
Set-Cell $ws "J" 79 3133.3333
Set-Cell $ws "K" 79 3260.8262
Set-Cell $ws "L" 79 3133.3333
Set-Cell $ws "M" 79 -2168.8262
Set-Cell $ws "N" 79 -5317.3333

Set-Cell $ws "H" 92 27780106
Set-Cell $ws "I" 92 33335376
Set-Cell $ws "J" 92 3750
Set-Cell $ws "K" 92 33335376
Set-Cell $ws "L" 92 3750
Set-Cell $ws "M" 92 -33334128

Set-Cell $ws "H" 120 37500
Set-Cell $ws "I" 120 0
Set-Cell $ws "J" 120 37500
Set-Cell $ws "K" 120 0
Set-Cell $ws "L" 120 37500
Set-Cell $ws "N" 120 -47176

Set-Cell $ws "H" 125 1478.9231
Set-Cell $ws "I" 125 788.6667
Set-Cell $ws "J" 125 2070.5715
Set-Cell $ws "K" 125 7098.0003
Set-Cell $ws "L" 125 18635.1435
Set-Cell $ws "M" 125 -4638.0003
Set-Cell $ws "N" 125 -23555.1435

Set-Cell $ws "H" 129 947.96875
Set-Cell $ws "I" 129 542.61536
Set-Cell $ws "J" 129 1011.4578
Set-Cell $ws "K" 129 1627.84608
Set-Cell $ws "L" 129 3034.3734
Set-Cell $ws "M" 129 3372.15392
Set-Cell $ws "N" 129 -13034.3734

Set-Cell $ws "H" 138 1830
Set-Cell $ws "I" 138 1649.8518
Set-Cell $ws "J" 138 2073.2
Set-Cell $ws "K" 138 4949.555399999999
Set-Cell $ws "L" 138 6219.599999999999
Set-Cell $ws "M" 138 190.4446000000007
Set-Cell $ws "N" 138 -16499.6


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
Set-Cell $ws "H" 32 409464.94
Set-Cell $ws "I" 32 457929.22
Set-Cell $ws "J" 32 74620.73
Set-Cell $ws "K" 32 457929.22
Set-Cell $ws "L" 32 74620.73
Set-Cell $ws "M" 32 -457642.22
Set-Cell $ws "N" 32 -75194.73

Set-Cell $ws "H" 45 3517.5715
Set-Cell $ws "I" 45 3251.7273
Set-Cell $ws "J" 45 3810
Set-Cell $ws "K" 45 3251.7273
Set-Cell $ws "L" 45 3810
Set-Cell $ws "M" 45 -2874.7273
Set-Cell $ws "N" 45 -4564

Set-Cell $ws "H" 63 5341.7144
Set-Cell $ws "I" 63 2971.1428
Set-Cell $ws "J" 63 7712.2856
Set-Cell $ws "K" 63 2971.1428
Set-Cell $ws "L" 63 7712.2856
Set-Cell $ws "M" 63 -2285.1428

Set-Cell $ws "H" 66 5341.7144
Set-Cell $ws "I" 66 2971.1428
Set-Cell $ws "J" 66 7712.2856
Set-Cell $ws "K" 66 14855.714
Set-Cell $ws "L" 66 38561.428
Set-Cell $ws "M" 66 -11423.714

Set-Cell $ws "H" 74 1327.8235
Set-Cell $ws "I" 74 941
Set-Cell $ws "J" 74 3133
Set-Cell $ws "K" 74 941
Set-Cell $ws "L" 74 3133
Set-Cell $ws "M" 74 -67

Set-Cell $ws "H" 77 1327.8235
Set-Cell $ws "I" 77 941
Set-Cell $ws "J" 77 3133
Set-Cell $ws "K" 77 4705
Set-Cell $ws "L" 77 15665
Set-Cell $ws "M" 77 -337

Set-Cell $ws "H" 132 3140.2246
Set-Cell $ws "I" 132 1810.2572
Set-Cell $ws "J" 132 6465.143
Set-Cell $ws "K" 132 5430.7716
Set-Cell $ws "L" 132 19395.429
Set-Cell $ws "M" 132 -2900.7716


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
Set-Cell $ws "H" 20 1658.2667
Set-Cell $ws "I" 20 1572.4642
Set-Cell $ws "J" 20 1799.5883
Set-Cell $ws "K" 20 1572.4642
Set-Cell $ws "L" 20 1799.5883
Set-Cell $ws "M" 20 -1325.4642
Set-Cell $ws "N" 20 -2293.5883

Set-Cell $ws "H" 86 111113496
Set-Cell $ws "I" 86 200002200
Set-Cell $ws "J" 86 2599.75
Set-Cell $ws "K" 86 200002200
Set-Cell $ws "L" 86 2599.75
Set-Cell $ws "M" 86 -200001077
Set-Cell $ws "N" 86 -4845.75

Set-Cell $ws "H" 89 111113496
Set-Cell $ws "I" 89 200002200
Set-Cell $ws "J" 89 2599.75
Set-Cell $ws "K" 89 1000011000
Set-Cell $ws "L" 89 12998.75
Set-Cell $ws "M" 89 -1000005384
Set-Cell $ws "N" 89 -24230.75

Set-Cell $ws "H" 94 1511.6842
Set-Cell $ws "I" 94 1232.625
Set-Cell $ws "J" 94 3000
Set-Cell $ws "K" 94 1232.625
Set-Cell $ws "L" 94 3000
Set-Cell $ws "M" 94 -781.625
Set-Cell $ws "N" 94 -3902


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
Set-Cell $ws "H" 63 0
Set-Cell $ws "I" 63 0
Set-Cell $ws "J" 63 0
Set-Cell $ws "K" 63 0
Set-Cell $ws "L" 63 0
$ws.Range("N63").ClearContents()

Set-Cell $ws "H" 66 0
Set-Cell $ws "I" 66 0
Set-Cell $ws "J" 66 0
Set-Cell $ws "K" 66 0
Set-Cell $ws "L" 66 0
$ws.Range("N66").ClearContents()


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
Set-Cell $ws "H" 5 326.8889
Set-Cell $ws "I" 5 301
Set-Cell $ws "J" 5 1000
Set-Cell $ws "K" 5 903
Set-Cell $ws "L" 5 3000
Set-Cell $ws "M" 5 -791

Set-Cell $ws "H" 75 648
Set-Cell $ws "I" 75 434
Set-Cell $ws "J" 75 933.3333
Set-Cell $ws "K" 75 1302
Set-Cell $ws "L" 75 2799.9999
Set-Cell $ws "M" 75 -304
Set-Cell $ws "N" 75 -4795.9999

Set-Cell $ws "H" 78 648
Set-Cell $ws "I" 78 434
Set-Cell $ws "J" 78 933.3333
Set-Cell $ws "K" 78 3906
Set-Cell $ws "L" 78 8399.9997
Set-Cell $ws "M" 78 1086
Set-Cell $ws "N" 78 -18383.9997

Set-Cell $ws "H" 117 449.16666
Set-Cell $ws "I" 117 504.5
Set-Cell $ws "J" 117 421.5
Set-Cell $ws "K" 117 1513.5
Set-Cell $ws "L" 117 1264.5
Set-Cell $ws "M" 117 1928.5
Set-Cell $ws "N" 117 -8148.5

Set-Cell $ws "H" 122 6867.6875
Set-Cell $ws "I" 122 529.61536
Set-Cell $ws "J" 122 34332.668
Set-Cell $ws "K" 122 4766.53824
Set-Cell $ws "L" 122 308994.012
Set-Cell $ws "M" 122 -2316.53824
Set-Cell $ws "N" 122 -313894.012

Set-Cell $ws "H" 131 1082.907
Set-Cell $ws "I" 131 915
Set-Cell $ws "J" 131 1091.0975
Set-Cell $ws "K" 131 2745
Set-Cell $ws "L" 131 3273.2925
Set-Cell $ws "M" 131 2295
Set-Cell $ws "N" 131 -13353.2925

Set-Cell $ws "H" 135 326.8889
Set-Cell $ws "I" 135 301
Set-Cell $ws "J" 135 1000
Set-Cell $ws "K" 135 2709
Set-Cell $ws "L" 135 9000
Set-Cell $ws "M" 135 -174

Set-Cell $ws "H" 139 3583.1875
Set-Cell $ws "I" 139 1742.2222
Set-Cell $ws "J" 139 4303.5654
Set-Cell $ws "K" 139 5226.6666
Set-Cell $ws "L" 139 12910.6962
Set-Cell $ws "M" 139 -86.66659999999956
Set-Cell $ws "N" 139 -23190.6962


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
Set-Cell $ws "H" 70 5652.5405
Set-Cell $ws "I" 70 5643.1665
Set-Cell $ws "J" 70 5669.846
Set-Cell $ws "K" 70 5643.1665
Set-Cell $ws "L" 70 5669.846
Set-Cell $ws "M" 70 -5373.1665
Set-Cell $ws "N" 70 -6209.846

Set-Cell $ws "H" 73 5652.5405
Set-Cell $ws "I" 73 5643.1665
Set-Cell $ws "J" 73 5669.846
Set-Cell $ws "K" 73 5643.1665
Set-Cell $ws "L" 73 5669.846
Set-Cell $ws "M" 73 -4707.1665
Set-Cell $ws "N" 73 -7541.846

Set-Cell $ws "H" 80 36388430
Set-Cell $ws "I" 80 50902904
Set-Cell $ws "J" 80 102250
Set-Cell $ws "K" 80 50902904
Set-Cell $ws "L" 80 102250
Set-Cell $ws "M" 80 -50901906
Set-Cell $ws "N" 80 -104246

Set-Cell $ws "H" 82 25999.445
Set-Cell $ws "I" 82 0
Set-Cell $ws "J" 82 25999.445
Set-Cell $ws "K" 82 0
Set-Cell $ws "L" 82 25999.445
Set-Cell $ws "N" 82 -26765.445

Set-Cell $ws "H" 83 36388430
Set-Cell $ws "I" 83 50902904
Set-Cell $ws "J" 83 102250
Set-Cell $ws "K" 83 254514520
Set-Cell $ws "L" 83 511250
Set-Cell $ws "M" 83 -254509528
Set-Cell $ws "N" 83 -521234

Set-Cell $ws "H" 85 25999.445
Set-Cell $ws "I" 85 0
Set-Cell $ws "J" 85 25999.445
Set-Cell $ws "K" 85 0
Set-Cell $ws "L" 85 25999.445
Set-Cell $ws "N" 85 -28651.445

Set-Cell $ws "H" 132 2536.8276
Set-Cell $ws "I" 132 2371.6667
Set-Cell $ws "J" 132 2807.0908
Set-Cell $ws "K" 132 7115.000100000001
Set-Cell $ws "L" 132 8421.2724
Set-Cell $ws "M" 132 -4585.000100000001
Set-Cell $ws "N" 132 -13481.2724


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
Set-Cell $ws "H" 68 1851.3077
Set-Cell $ws "I" 68 1760.7273
Set-Cell $ws "J" 68 2349.5
Set-Cell $ws "K" 68 1760.7273
Set-Cell $ws "L" 68 2349.5
Set-Cell $ws "M" 68 -1011.7273
Set-Cell $ws "N" 68 -3847.5

Set-Cell $ws "H" 71 1851.3077
Set-Cell $ws "I" 71 1760.7273
Set-Cell $ws "J" 71 2349.5
Set-Cell $ws "K" 71 8803.636500000001
Set-Cell $ws "L" 71 11747.5
Set-Cell $ws "M" 71 -5059.636500000001
Set-Cell $ws "N" 71 -19235.5

Set-Cell $ws "H" 82 33336724
Set-Cell $ws "I" 82 41669956
Set-Cell $ws "J" 82 3800
Set-Cell $ws "K" 82 41669956
Set-Cell $ws "L" 82 3800
Set-Cell $ws "M" 82 -41669595
Set-Cell $ws "N" 82 -4522

Set-Cell $ws "H" 85 33336724
Set-Cell $ws "I" 85 41669956
Set-Cell $ws "J" 85 3800
Set-Cell $ws "K" 85 41669956
Set-Cell $ws "L" 85 3800
Set-Cell $ws "M" 85 -41668708
Set-Cell $ws "N" 85 -6296

Set-Cell $ws "H" 93 9184
Set-Cell $ws "I" 93 12588
Set-Cell $ws "J" 93 1525
Set-Cell $ws "K" 93 12588
Set-Cell $ws "L" 93 1525
Set-Cell $ws "M" 93 -11340
Set-Cell $ws "N" 93 -4021

Set-Cell $ws "H" 132 3400.4
Set-Cell $ws "I" 132 2926.1667
Set-Cell $ws "J" 132 4111.75
Set-Cell $ws "K" 132 8778.500100000001
Set-Cell $ws "L" 132 12335.25
Set-Cell $ws "M" 132 -6248.500100000001
Set-Cell $ws "N" 132 -17395.25

